# Daily attendance processing - 2025-11-18 17:21:51
# Normalize the "Recorded By" (column G) values: when the value starts with
# "System, " (i.e. "System" listed first among multiple recorders), move
# "System" to the end of the comma-separated list instead of the front.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val.StartsWith("System, ")) {
        $rest = $val.Substring(8)
        $cell.Value2 = "$rest, System"
    }
}
